$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1 and fill in the Save values for each data row
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1

# Match H1's formatting to the other header cells (e.g. G1 "sum") - bold,
# bordered, centered header style
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
